$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 19:22"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 488755
$ws.Range("C4").Value = 20189
$ws.Range("D4").Value = 26163
$ws.Range("E4").Value = 444597
$ws.Range("F4").Value = 10896
$ws.Range("G4").Value = 1304
$ws.Range("H4").Value = 17995

# Row 10: Reino Unido
$ws.Range("A10").Value = "Reino Unido"
$ws.Range("B10").Value = 73758
$ws.Range("C10").Value = 8681
$ws.Range("D10").Value = 135
$ws.Range("E10").Value = 64665
$ws.Range("F10").Value = 1559
$ws.Range("G10").Value = 980
$ws.Range("H10").Value = 8958

# Row 19: Austria
$ws.Range("A19").Value = "Austria"
$ws.Range("B19").Value = 13549
$ws.Range("C19").Value = 305
$ws.Range("D19").Value = 6064
$ws.Range("E19").Value = 7166
$ws.Range("F19").Value = 261
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 319

# Row 28: Noruega
$ws.Range("A28").Value = "Noruega"
$ws.Range("B28").Value = 6298
$ws.Range("C28").Value = 79
$ws.Range("D28").Value = 32
$ws.Range("E28").Value = 6154
$ws.Range("F28").Value = 70
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 112

# Row 47: Republica Dominicana
$ws.Range("A47").Value = "Republica Dominicana"
$ws.Range("B47").Value = 2620
$ws.Range("C47").Value = 271
$ws.Range("D47").Value = 98
$ws.Range("E47").Value = 2396
$ws.Range("F47").Value = 147
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 126

# Row 48: Catar
$ws.Range("A48").Value = "Catar"
$ws.Range("B48").Value = 2512
$ws.Range("C48").Value = 136
$ws.Range("D48").Value = 227
$ws.Range("E48").Value = 2279
$ws.Range("F48").Value = 37
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 6

# Row 49: Tailandia
$ws.Range("A49").Value = "Tailandia"
$ws.Range("B49").Value = 2473
$ws.Range("C49").Value = 50
$ws.Range("D49").Value = 1013
$ws.Range("E49").Value = 1427
$ws.Range("F49").Value = 61
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 33

# Row 53: Grecia
$ws.Range("A53").Value = "Grecia"
$ws.Range("B53").Value = 2011
$ws.Range("C53").Value = 56
$ws.Range("D53").Value = 269
$ws.Range("E53").Value = 1651
$ws.Range("F53").Value = 77
$ws.Range("G53").Value = 4
$ws.Range("H53").Value = 91

# Row 56: Argentina
$ws.Range("A56").Value = "Argentina"
$ws.Range("B56").Value = 1894
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 375
$ws.Range("E56").Value = 1438
$ws.Range("F56").Value = 96
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 81

# Row 61: Marruecos
$ws.Range("A61").Value = "Marruecos"
$ws.Range("B61").Value = 1448
$ws.Range("C61").Value = 74
$ws.Range("D61").Value = 122
$ws.Range("E61").Value = 1219
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 10
$ws.Range("H61").Value = 107

# Row 62: Moldavia
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 1438
$ws.Range("C62").Value = 149
$ws.Range("D62").Value = 56
$ws.Range("E62").Value = 1353
$ws.Range("F62").Value = 80
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 29

# Row 108: Estado de Palestina
$ws.Range("A108").Value = "Estado de Palestina"
$ws.Range("B108").Value = 267
$ws.Range("C108").Value = 4
$ws.Range("D108").Value = 45
$ws.Range("E108").Value = 220
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 2

# Row 134: Aruba
$ws.Range("A134").Value = "Aruba"
$ws.Range("B134").Value = 86
$ws.Range("C134").Value = 4
$ws.Range("D134").Value = 27
$ws.Range("E134").Value = 59
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

# Row 135: Monaco
$ws.Range("A135").Value = "Monaco"
$ws.Range("B135").Value = 84
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 5
$ws.Range("E135").Value = 78
$ws.Range("F135").Value = 4
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 1

# Row 136: Guayana Francesa
$ws.Range("A136").Value = "Guayana Francesa"
$ws.Range("B136").Value = 83
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 43
$ws.Range("E136").Value = 40
$ws.Range("F136").Value = 1
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

# Row 164: Somalia
$ws.Range("A164").Value = "Somalia"
$ws.Range("B164").Value = 21
$ws.Range("C164").Value = 9
$ws.Range("D164").Value = 1
$ws.Range("E164").Value = 19
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1

# Row 165: Mozambique
$ws.Range("A165").Value = "Mozambique"
$ws.Range("B165").Value = 20
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 2
$ws.Range("E165").Value = 18
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0

# Row 166: Antigua y Barbuda
$ws.Range("A166").Value = "Antigua y Barbuda"
$ws.Range("B166").Value = 19
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 17
$ws.Range("F166").Value = 1
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 2

# Row 167: Angola
$ws.Range("A167").Value = "Angola"
$ws.Range("B167").Value = 19
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 2
$ws.Range("E167").Value = 15
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 2

# Row 168: Siria
$ws.Range("A168").Value = "Siria"
$ws.Range("B168").Value = 19
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 4
$ws.Range("E168").Value = 13
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 2

# Row 169: Maldivas
$ws.Range("A169").Value = "Maldivas"
$ws.Range("B169").Value = 19
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 13
$ws.Range("E169").Value = 6
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0

# Row 170: Nueva Caledonia
$ws.Range("A170").Value = "Nueva Caledonia"
$ws.Range("B170").Value = 18
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 1
$ws.Range("E170").Value = 17
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

# Row 171: Guinea Ecuatorial
$ws.Range("A171").Value = "Guinea Ecuatorial"
$ws.Range("B171").Value = 18
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 3
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

# Row 172: Islas Virgenes de los Estados Unidos
$ws.Range("A172").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B172").Value = 17
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 17
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

# Row 173: Laos
$ws.Range("A173").Value = "Laos"
$ws.Range("B173").Value = 16
$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 16
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

# Row 174: Fiyi
$ws.Range("A174").Value = "Fiyi"
$ws.Range("B174").Value = 16
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 16
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# Row 175: Namibia
$ws.Range("A175").Value = "Namibia"
$ws.Range("B175").Value = 16
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 3
$ws.Range("E175").Value = 13
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

# Row 176: Mongolia
$ws.Range("A176").Value = "Mongolia"
$ws.Range("B176").Value = 16
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 4
$ws.Range("E176").Value = 12
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

# Row 177: Dominica
$ws.Range("A177").Value = "Dominica"
$ws.Range("B177").Value = 16
$ws.Range("C177").Value = 1
$ws.Range("D177").Value = 5
$ws.Range("E177").Value = 11
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

# Row 178: Sudan
$ws.Range("A178").Value = "Sudan"
$ws.Range("B178").Value = 15
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 2
$ws.Range("E178").Value = 11
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 2

# Row 179: Santa Lucia
$ws.Range("A179").Value = "Santa Lucia"
$ws.Range("B179").Value = 14
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 1
$ws.Range("E179").Value = 13
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

# Row 180: Curazao
$ws.Range("A180").Value = "Curazao"
$ws.Range("B180").Value = 14
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 7
$ws.Range("E180").Value = 6
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 1

# Row 181: Botsuana
$ws.Range("A181").Value = "Botsuana"
$ws.Range("B181").Value = 13
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 12
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 1

# Row 182: Granada
$ws.Range("A182").Value = "Granada"
$ws.Range("B182").Value = 12
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 12
$ws.Range("F182").Value = 2
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0

# Row 183: San Vicente y las Granadinas
$ws.Range("A183").Value = "San Vicente y las Granadinas"
$ws.Range("B183").Value = 12
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 1
$ws.Range("E183").Value = 11
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

